$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated s_vals data (regen to filter save games)
$data = @{
    2 = @{ B = 1.505614041169197;  C = 1.65323645889881;  D = 0.7127328510149897; E = 0.4998867070740569; G = 4.371470058157054 }
    3 = @{ B = 1.505614041169197;  C = 1.65323645889881;  D = 0.1529057820181812; E = 0.4998867070740569; G = 3.811642989160245 }
    4 = @{ B = 0.06328177979961902; C = 0.3375848360084654; D = 0.7127328510149897; E = 0.4998867070740569; G = 1.613486173897131 }
    5 = @{ B = 0.7287194209349384;  C = 1.65323645889881;  D = 0.1529057820181812; E = 0.4998867070740569; G = 3.034748368925986 }
    6 = @{ B = 3.182878228561681;   C = 1.65323645889881;  D = 0.7127328510149897; E = 0.4998867070740569; G = 6.048734245549538 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
